$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.632.22"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.544.18"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'317.93"
$ws.Range("E5").Value = "  +4.85%  "
$ws.Range("D6").Value = "'94.99"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "'36.49"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'7.67"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.114"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "2.934.01"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "'15.79"
$ws.Range("E15").Value = "  +5.42%  "
$ws.Range("D16").Value = "2.541.58"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'0.868"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "42.682.45"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'6.64"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "0.0₃0968"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "'71.03"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "'253.10"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "'27.35"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +4.64%  "
$ws.Range("D29").Value = "'10.23"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "'39.28"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'155.33"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "'2.17"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("D35").Value = "'19.34"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'0.0789"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "'23.90"
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("D41").Value = "'2.36"
$ws.Range("E41").Value = "  +9.48%  "
$ws.Range("D42").Value = "'3.84"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'3.35"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "'0.0303"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "2.036.17"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "'84.70"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("D48").Value = "'8.95"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "2.788.32"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'74.10"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -0.35%  "
